$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 21) with the new problem entry
# (set B before A so the shared-strings table order matches: class name first, description second)
$ws.Range("B21").Value = "RearrangeAlternatingPositiveAndNegativeNumbers"
$ws.Range("A21").Value = "Rearrange array in alternating positive & negative items with O(1) extra space"

# Update selection to match the new active cell
$ws.Range("A21").Select()
